# Updated symbols/footprint for USB and 5K Trimmer Pots.
#
# The original "cpu" sheet is renamed to "_cpu" (kept intact as-is), and a
# new "cpu" sheet is inserted right after it. The new sheet starts out with
# just the BOM header row (plus three new columns: Vendor\Part #, Vendor
# URL, Description) and becomes the active/selected tab.

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet, and remember its current selection -------
$wsOld = $wb.Worksheets.Item(1)
$wsOld.Name = "_cpu"
[void]$wsOld.Columns.Item(6).Select()

# --- Insert the new "cpu" sheet right after the old one -------------------
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsOld)
$wsNew.Name = "cpu"

# --- Copy the header-row formatting (bold/fill) from the old sheet --------
$wsOld.Range("A1:H1").Copy()
$wsNew.Range("A1:I1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Header row -------------------------------------------------------
$wsNew.Range("A1").Value = "#"
$wsNew.Range("B1").Value = "Reference"
$wsNew.Range("C1").Value = "Qty"
$wsNew.Range("D1").Value = "Value"
$wsNew.Range("E1").Value = "Footprint"
$wsNew.Range("F1").Value = "Vendor\Part #"
$wsNew.Range("G1").Value = "Vendor URL"
$wsNew.Range("H1").Value = "Datasheet"
$wsNew.Range("I1").Value = "Description"

# --- Column widths ------------------------------------------------------
$wsNew.Columns.Item(2).ColumnWidth = 18.666666666666668
$wsNew.Columns.Item(4).ColumnWidth = 25
$wsNew.Columns.Item(5).ColumnWidth = 35.833333333333336
$wsNew.Columns.Item(6).ColumnWidth = 36.833333333333336
$wsNew.Columns.Item(7).ColumnWidth = 30.666666666666668
$wsNew.Columns.Item(8).ColumnWidth = 31.166666666666668
$wsNew.Columns.Item(9).ColumnWidth = 25.833333333333332

# --- Make the new sheet active / selected --------------------------------
[void]$wsNew.Range("A2").Select()
[void]$wsNew.Activate()
